# The edit reorders the data rows (2-20) of the "Artfynd" sheet: every
# row's full content (all columns A:AY) is an exact match of some other
# row's content from before the edit - i.e. the rows were shuffled
# around while row 1 (header) and row 5 stay in place.
#
# Mapping: destination row number -> source row number (content that
# should end up there), derived by matching each row's full content
# before/after the edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    2  = 7
    3  = 8
    4  = 13
    6  = 14
    7  = 16
    8  = 3
    9  = 15
    10 = 11
    11 = 4
    12 = 2
    13 = 20
    14 = 17
    15 = 19
    16 = 9
    17 = 18
    18 = 12
    19 = 6
    20 = 10
}

$lastCol = "AY"
$stagingOffset = 1000

# Record, for every row about to be used as a source, whether the
# "sometimes populated" columns (K/L/M/N/AC) actually hold data in the
# original layout - Range.Copy() stamps blank placeholder cells across
# the whole A:AY block it touches, so after the moves we need to know
# which of those columns to strip back out again for rows whose source
# content didn't originally have them. K/L/M/N are always blank even on
# rows that "have" the block, so use AC (which carries real text, e.g.
# "ringhack") as the tell-tale column instead. This probe runs first,
# while rows 2-20 still hold their original content, so the moves below
# can't taint the reading.
$sparseCols = @("K", "L", "M", "N", "AC")
$sourceHasSparse = @{}
foreach ($srcRow in $mapping.Values) {
    if (-not $sourceHasSparse.ContainsKey($srcRow)) {
        $probe = $ws.Range("AC" + $srcRow).Value2
        $sourceHasSparse[$srcRow] = ($probe -ne $null) -and ($probe -ne "")
    }
}

# Step 1: stage every affected row into a scratch area far below the
# data so that overlapping reads/writes during the permutation can't
# clobber a row before it's been copied elsewhere (the permutation has
# multiple independent cycles).
foreach ($srcRow in $mapping.Values) {
    $src = $ws.Range("A" + $srcRow + ":" + $lastCol + $srcRow)
    $stageRow = $stagingOffset + $srcRow
    $dst = $ws.Range("A" + $stageRow)
    $src.Copy($dst)
}

# Step 2: copy staged rows into their final destinations.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $stageRow = $stagingOffset + $srcRow
    $src = $ws.Range("A" + $stageRow + ":" + $lastCol + $stageRow)
    $dst = $ws.Range("A" + $destRow)
    $src.Copy($dst)
}

# Step 3: wipe the scratch area used for staging.
$ws.Range("A1002:" + $lastCol + "1020").Clear()

# Step 4: the whole-row Copy() calls above stamp blank placeholder
# cells across every column in A:AY, including columns that the pasted
# source row never actually used. Clear those back out so the sheet's
# cell footprint matches the source data exactly.
$neverUsedCols = @("J", "O", "X", "AF", "AH", "AI", "AJ", "AK", "AL", "AM", "AN", "AO", "AP", "AQ", "AR", "AS", "AU", "AV")
foreach ($col in $neverUsedCols) {
    $ws.Range($col + "2:" + $col + "20").Clear()
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    if (-not $sourceHasSparse[$srcRow]) {
        foreach ($col in $sparseCols) {
            $ws.Range($col + $destRow).Clear()
        }
    }
}
